$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update status text wherever it appears ("Ready for handoff" -> "Handback transform failed")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Widen the "Error Detail" column (column P, 16th column) on both language sheets
$wsZhCn.Columns.Item(16).ColumnWidth = 39.15
$wsDeDe.Columns.Item(16).ColumnWidth = 39.15

# Populate Error Detail (P3) with the handback transform failure message for each language
$wsZhCn.Range("P3").Value = "Handback file name: 4qh4cxv3.1xm is different with handoff file name: d0e45ec0-e79a-4e91-94c7-b7fc792b0187.3d319a04d0f6f78a2f8cd3bb0cb7702685af0477.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: 4qh4cxv3.1xm is different with handoff file name: d0e45ec0-e79a-4e91-94c7-b7fc792b0187.3d319a04d0f6f78a2f8cd3bb0cb7702685af0477.de-de."
